$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44558
$ws.Range("J2").Value = 250
$ws.Range("K2").Value = 15000
$ws.Range("L2").Value = 16000
$ws.Range("M2").Value = 15400
$ws.Range("O2").Value = "Región Metropolitana"
$ws.Range("P2").Value = 616

# Row 3
$ws.Range("D3").Value = 44251
$ws.Range("K3").Value = 27000
$ws.Range("L3").Value = 28000
$ws.Range("M3").Value = 27500
$ws.Range("O3").Value = "Región Metropolitana"
$ws.Range("P3").Value = 1100

# Row 4
$ws.Range("D4").Value = 44308
$ws.Range("K4").Value = 28000
$ws.Range("L4").Value = 30000
$ws.Range("M4").Value = 29000
$ws.Range("P4").Value = 1160

# Row 5
$ws.Range("D5").Value = 44216
$ws.Range("K5").Value = 26000
$ws.Range("L5").Value = 28000
$ws.Range("M5").Value = 27000
$ws.Range("O5").Value = "Región del Maule"
$ws.Range("P5").Value = 1080

# Row 6
$ws.Range("D6").Value = 44342
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 28000
$ws.Range("L6").Value = 30000
$ws.Range("M6").Value = 29000
$ws.Range("O6").Value = "Región Metropolitana"
$ws.Range("P6").Value = 1160

# Row 7
$ws.Range("D7").Value = 44580
$ws.Range("K7").Value = 28000
$ws.Range("L7").Value = 30000
$ws.Range("M7").Value = 29000
$ws.Range("O7").Value = "Región Metropolitana"
$ws.Range("P7").Value = 1160

# Row 8
$ws.Range("D8").Value = 44587
$ws.Range("J8").Value = 220
$ws.Range("K8").Value = 23000
$ws.Range("L8").Value = 24000
$ws.Range("M8").Value = 23545
$ws.Range("P8").Value = 942

# Row 9
$ws.Range("D9").Value = 44594
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = 24000
$ws.Range("L9").Value = 25000
$ws.Range("M9").Value = 24500
$ws.Range("P9").Value = 980

# Row 10
$ws.Range("D10").Value = 44203
$ws.Range("O10").Value = "Región de O'Higgins"

# Row 11
$ws.Range("D11").Value = 44574
$ws.Range("K11").Value = 30000
$ws.Range("L11").Value = 32000
$ws.Range("M11").Value = 31000
$ws.Range("P11").Value = 1240

# Row 12
$ws.Range("D12").Value = 44328
$ws.Range("K12").Value = 32000
$ws.Range("L12").Value = 34000
$ws.Range("M12").Value = 33000
$ws.Range("O12").Value = "Región Metropolitana"
$ws.Range("P12").Value = 1320

# Row 14
$ws.Range("D14").Value = 44265
$ws.Range("K14").Value = 22000
$ws.Range("L14").Value = 24000
$ws.Range("M14").Value = 23000
$ws.Range("P14").Value = 920

# Row 15
$ws.Range("D15").Value = 44316
$ws.Range("K15").Value = 26000
$ws.Range("L15").Value = 27000
$ws.Range("M15").Value = 26500
$ws.Range("O15").Value = "Región Metropolitana"
$ws.Range("P15").Value = 1060

# Row 16
$ws.Range("D16").Value = 44320
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 26000
$ws.Range("L16").Value = 28000
$ws.Range("M16").Value = 27000
$ws.Range("O16").Value = "Región del Maule"
$ws.Range("P16").Value = 1080

# Row 17
$ws.Range("D17").Value = 44560
$ws.Range("K17").Value = 25000
$ws.Range("L17").Value = 26000
$ws.Range("M17").Value = 25500
$ws.Range("O17").Value = "Región del Maule"
$ws.Range("P17").Value = 1020

# Row 18
$ws.Range("D18").Value = 44313
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = 30000
$ws.Range("L18").Value = 32000
$ws.Range("M18").Value = 31000
$ws.Range("P18").Value = 1240

# Row 19
$ws.Range("D19").Value = 44244
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = 25000
$ws.Range("L19").Value = 26000
$ws.Range("M19").Value = 25500
$ws.Range("O19").Value = "Región del Maule"
$ws.Range("P19").Value = 1020

# Row 20
$ws.Range("D20").Value = 44194
$ws.Range("K20").Value = 30000
$ws.Range("L20").Value = 32000
$ws.Range("M20").Value = 31000
$ws.Range("P20").Value = 1240

# Row 21
$ws.Range("D21").Value = 44210
$ws.Range("K21").Value = 32000
$ws.Range("L21").Value = 34000
$ws.Range("M21").Value = 33000
$ws.Range("O21").Value = "Región del Maule"
$ws.Range("P21").Value = 1320

# Row 22
$ws.Range("D22").Value = 44188
$ws.Range("K22").Value = 42000
$ws.Range("L22").Value = 44000
$ws.Range("M22").Value = 43000
$ws.Range("O22").Value = "Región de O'Higgins"
$ws.Range("P22").Value = 1720

# Row 23
$ws.Range("D23").Value = 44568
$ws.Range("J23").Value = 200
$ws.Range("K23").Value = 25000
$ws.Range("L23").Value = 26000
$ws.Range("M23").Value = 25500
$ws.Range("O23").Value = "Región de O'Higgins"
$ws.Range("P23").Value = 1020

# Row 24
$ws.Range("D24").Value = 44294
$ws.Range("K24").Value = 20000
$ws.Range("L24").Value = 22000
$ws.Range("M24").Value = 21000
$ws.Range("O24").Value = "Región del Maule"
$ws.Range("P24").Value = 840

# Row 25
$ws.Range("D25").Value = 44279
$ws.Range("K25").Value = 28000
$ws.Range("L25").Value = 30000
$ws.Range("M25").Value = 29000
$ws.Range("O25").Value = "Región del Maule"
$ws.Range("P25").Value = 1160

# Row 26
$ws.Range("D26").Value = 44236
$ws.Range("K26").Value = 25000
$ws.Range("L26").Value = 26000
$ws.Range("M26").Value = 25500
$ws.Range("O26").Value = "Región del Maule"
$ws.Range("P26").Value = 1020

# Row 27
$ws.Range("D27").Value = 44272
$ws.Range("K27").Value = 22000
$ws.Range("L27").Value = 24000
$ws.Range("M27").Value = 23000
$ws.Range("P27").Value = 920

# Row 28
$ws.Range("D28").Value = 44349
$ws.Range("J28").Value = 60
$ws.Range("K28").Value = 30000
$ws.Range("L28").Value = 32000
$ws.Range("M28").Value = 31000
$ws.Range("O28").Value = "Región Metropolitana"
$ws.Range("P28").Value = 1240

# New row 29
$ws.Range("A29").Value = 11
$ws.Range("B29").Value = "Vega Monumental Concepción"
$ws.Range("C29").Value = "Bíobío"
$ws.Range("D29").Value = 44552
$ws.Range("E29").Value = 8
$ws.Range("F29").Value = 100112030
$ws.Range("G29").Value = "Poroto granado"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 100
$ws.Range("K29").Value = 38000
$ws.Range("L29").Value = 40000
$ws.Range("M29").Value = 39000
$ws.Range("N29").Value = "`$/saco 25 kilos"
$ws.Range("O29").Value = "Provincia de Limarí"
$ws.Range("P29").Value = 1560
$ws.Range("Q29").Value = 25
$ws.Range("R29").Value = "Hortaliza"
$ws.Range("D29").NumberFormat = "YYYY-MM-DD HH:MM:SS"